# TLRatio.xlsx update
# - Biomass calculation columns (F) converted from live formulas to static
#   computed values, with a few new helper cells (G) added for boxplot/PCO
#   notes.
# - Five new species rows appended (Caranx lugubris, Canthidermis maculata,
#   Carcharhinus falciformis, Thunnus albacares, Euthynnus lineatus).
# - TLRatio number format tightened from 3 to 5 decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Freeze the TLRatio (column F) formulas into plain computed values and
#    tighten the number format used for that ratio (0.000 -> 0.00000).
#    A blank "notes" cell (column G) is added alongside several of them.
# ---------------------------------------------------------------------
$ws.Range("F11").Value = 0.95877

$ws.Range("F12").Value = 0.97561
$ws.Range("G12").NumberFormat = "0.00000"

$ws.Range("F13").Value = 0.7755
$ws.Range("G13").NumberFormat = "0.00000"

$ws.Range("F14").Value = 0.93633
$ws.Range("G14").NumberFormat = "0.00000"

$ws.Range("F15").Value = 0.90334
$ws.Range("G15").NumberFormat = "0.00000"

$ws.Range("F16").Value = 0.82237
$ws.Range("G16").NumberFormat = "0.00000"

$ws.Range("F17").Value = 0.88067
$ws.Range("G17").NumberFormat = "0.00000"

$ws.Range("F11:F17").NumberFormat = "0.00000"

# ---------------------------------------------------------------------
# 2. Append the newly added species rows (25-29).
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Caranx lugubris"
$ws.Range("B25").Value = "Caranx lugubris"
$ws.Range("C25").Value = 0.0251
$ws.Range("D25").Value = 2.84
$ws.Range("E25").Value = "FL"
$ws.Range("F25").Value = 0.863

$ws.Range("A26").Value = "Canthidermis maculata"
$ws.Range("B26").Value = "Canthidermis maculata"
$ws.Range("C26").Value = 0.0176
$ws.Range("D26").Value = 3.055
$ws.Range("E26").Value = "FL"
$ws.Range("F26").Value = 1

$ws.Range("A27").Value = "Carcharhinus falciformis"
$ws.Range("B27").Value = "Carcharhinus falciformis"
$ws.Range("C27").Value = 0.0464
$ws.Range("D27").Value = 2.75
$ws.Range("E27").Value = "SL"
$ws.Range("F27").Value = 0.81722
$ws.Range("G27").NumberFormat = "0.00000"

$ws.Range("A28").Value = "Thunnus albacares"
$ws.Range("B28").Value = "Thunnus albacares"
$ws.Range("C28").Value = 0.0216
$ws.Range("D28").Value = 2.981
$ws.Range("E28").Value = "TL"
$ws.Range("F28").Value = 1

$ws.Range("A29").Value = "Euthynnus lineatus"
$ws.Range("B29").Value = "Euthynnus lineatus"
$ws.Range("C29").Value = 0.01
$ws.Range("D29").Value = 3.05
$ws.Range("E29").Value = "TL"
$ws.Range("F29").Value = 1

# ---------------------------------------------------------------------
# 3. New column widths for the columns (G/H/I) introduced to the right of
#    TLRatio.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10.5546875
$ws.Columns.Item(8).ColumnWidth = 10.21875
$ws.Columns.Item(9).ColumnWidth = 10.77734375

# ---------------------------------------------------------------------
# 4. Scroll/selection state: the sheet was left scrolled down a bit with
#    A30 (first blank row after the new data) selected.
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("A30").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
